$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "exclude" flags in column C
$ws.Range("C9").Value = 0
$ws.Range("C11").Value = 1

# Update the active selection on the sheet
$ws.Range("K9").Select()
